$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $ws.Range($cell).Style
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $text
    $ws.Range($cell).Style = $origStyle
}

$ws.Range("D2").Value = "23.413.79"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "1.643.38"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.11%  "
Set-TextValue "D5" "1.001"
$ws.Range("E5").Value = "  +0.03%  "
Set-TextValue "D6" "299.05"
$ws.Range("E6").Value = "  -1.68%  "
Set-TextValue "D7" "0.3785"
$ws.Range("E7").Value = "  -0.60%  "
Set-TextValue "D8" "0.3501"
$ws.Range("E8").Value = "  -3.68%  "
Set-TextValue "D9" "49.89"
$ws.Range("E9").Value = "  -2.09%  "
Set-TextValue "D10" "0.08070"
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("E11").Value = "  -3.89%  "
Set-TextValue "D13" "22.05"
Set-TextValue "D14" "6.362"
$ws.Range("E14").Value = "  -2.94%  "
Set-TextValue "D15" "7.312"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("E16").Value = "  -3.49%  "
$ws.Range("D17").Value = "1.642.50"
$ws.Range("E17").Value = "  -0.80%  "
Set-TextValue "D18" "96.46"
$ws.Range("E18").Value = "  -1.44%  "
Set-TextValue "D19" "0.07016"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("E21").Value = "  -2.78%  "
$ws.Range("E22").Value = "  -0.02%  "
Set-TextValue "D23" "12.33"
$ws.Range("E23").Value = "  -4.25%  "
$ws.Range("D24").Value = "23.447.01"
$ws.Range("E24").Value = "  -1.41%  "
Set-TextValue "D25" "2.491"
$ws.Range("E25").Value = "  -2.52%  "
Set-TextValue "D26" "2.912"
$ws.Range("E26").Value = "  -5.67%  "
Set-TextValue "D27" "20.81"
$ws.Range("E27").Value = "  -2.63%  "
Set-TextValue "D28" "152.65"
$ws.Range("E28").Value = "  +0.64%  "
Set-TextValue "D29" "5.204"
$ws.Range("E29").Value = "  -0.60%  "
Set-TextValue "D30" "132.56"
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("D31").Value = "1.825.66"
$ws.Range("E31").Value = "  -0.76%  "
Set-TextValue "D33" "2.117"
$ws.Range("E33").Value = "  -3.48%  "
$ws.Range("E34").Value = "  -4.15%  "
Set-TextValue "D35" "0.9783"
$ws.Range("E35").Value = "  -9.53%  "
Set-TextValue "D36" "0.02699"
$ws.Range("E36").Value = "  -4.96%  "
Set-TextValue "D37" "0.08733"
$ws.Range("E37").Value = "  -1.13%  "
Set-TextValue "D38" "0.2427"
$ws.Range("E38").Value = "  -4.22%  "
Set-TextValue "D39" "5.900"
$ws.Range("E39").Value = "  -4.09%  "
Set-TextValue "D40" "0.06803"
$ws.Range("E40").Value = "  -4.58%  "
Set-TextValue "D41" "12.82"
$ws.Range("E41").Value = "  -3.48%  "
Set-TextValue "D42" "0.6861"
$ws.Range("E42").Value = "  -3.27%  "
Set-TextValue "D43" "1.299"
$ws.Range("E43").Value = "  -3.61%  "
Set-TextValue "D44" "15.67"
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("E45").Value = "  +0.06%  "
Set-TextValue "D46" "0.6332"
$ws.Range("E46").Value = "  -3.71%  "
Set-TextValue "D47" "2.247"
$ws.Range("E47").Value = "  -3.94%  "
Set-TextValue "D48" "3.903"
$ws.Range("E48").Value = "  -1.62%  "
Set-TextValue "D49" "0.07710"
$ws.Range("E49").Value = "  -3.22%  "
Set-TextValue "D50" "126.99"
$ws.Range("E50").Value = "  -1.13%  "
Set-TextValue "D51" "1.141"
$ws.Range("E51").Value = "  -4.88%  "
